$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the cells we are about to update as Text so that numeric-looking
# strings (with significant trailing zeros, leading zeros, etc.) are
# preserved exactly as literal text instead of being auto-converted to numbers.
$updateCells = @("D2","E2","D3","E3","D4","E4","D5","E5","D6","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","D18","E18","D19","E19","E21","D22","E22","D23","E23","D25","E25","D26","E26","D27","E27","E28","D40","E40","D41","E41","D42","E42","D43","E43","D44","E44","D45","E45","D46","E46","D47","E47")
foreach ($addr in $updateCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "275.81"
$ws.Range("E2").Value = "0.69%"
$ws.Range("D3").Value = "27.19"
$ws.Range("E3").Value = "1.64%"
$ws.Range("D4").Value = "4.853"
$ws.Range("E4").Value = "-0.28%"
$ws.Range("D5").Value = "0.06393"
$ws.Range("E5").Value = "1.20%"
$ws.Range("D6").Value = "6.955"
$ws.Range("E6").Value = "1.04%"
$ws.Range("D7").Value = "1.194"
$ws.Range("E7").Value = "-7.13%"
$ws.Range("D8").Value = "0.8760"
$ws.Range("E8").Value = "0.57%"
$ws.Range("D9").Value = "0.1510"
$ws.Range("E9").Value = "3.46%"
$ws.Range("D10").Value = "0.05123"
$ws.Range("E10").Value = "1.40%"
$ws.Range("D11").Value = "0.07564"
$ws.Range("E11").Value = "2.23%"
$ws.Range("D12").Value = "0.02957"
$ws.Range("E12").Value = "-1.09%"
$ws.Range("D13").Value = "0.08983"
$ws.Range("E13").Value = "-0.58%"
$ws.Range("D14").Value = "0.001563"
$ws.Range("E14").Value = "-0.50%"
$ws.Range("D15").Value = "0.0006383"
$ws.Range("E15").Value = "1.58%"
$ws.Range("D16").Value = "0.006183"
$ws.Range("E16").Value = "4.98%"
$ws.Range("D17").Value = "3.472"
$ws.Range("E17").Value = "0.56%"
$ws.Range("D18").Value = "3.308"
$ws.Range("E18").Value = "-0.26%"
$ws.Range("D19").Value = "2.253"
$ws.Range("E19").Value = "-1.36%"
$ws.Range("E21").Value = "1.88%"
$ws.Range("D22").Value = "3.910"
$ws.Range("E22").Value = "-0.09%"
$ws.Range("D23").Value = "0.04411"
$ws.Range("E23").Value = "1.06%"
$ws.Range("D25").Value = "0.001176"
$ws.Range("E25").Value = "-0.11%"
$ws.Range("D26").Value = "0.003859"
$ws.Range("E26").Value = "-9.49%"
$ws.Range("D27").Value = "0.0001200"
$ws.Range("E27").Value = "0.02%"
$ws.Range("E28").Value = "14.68%"
$ws.Range("D40").Value = "0.04137"
$ws.Range("E40").Value = "2.51%"
$ws.Range("D41").Value = "0.006820"
$ws.Range("E41").Value = "2.22%"
$ws.Range("D42").Value = "0.1174"
$ws.Range("E42").Value = "0.58%"
$ws.Range("D43").Value = "0.002159"
$ws.Range("E43").Value = "3.37%"
$ws.Range("D44").Value = "0.01179"
$ws.Range("E44").Value = "-2.99%"
$ws.Range("D45").Value = "0.00005189"
$ws.Range("E45").Value = "-2.43%"
$ws.Range("D46").Value = "1.687"
$ws.Range("E46").Value = "-29.12%"
$ws.Range("D47").Value = "0.01852"
$ws.Range("E47").Value = "-7.36%"
